$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: insert the 3 new "full range" rows (p99-p100, p99.9-p100, p99.99-p100) ---
# Insert blank rows at (current) row 100, 110, 120 -- processed top-down so later
# positions already reflect the shift caused by earlier inserts.
$ws1.Rows("100:100").Insert() | Out-Null
$ws1.Rows("110:110").Insert() | Out-Null
$ws1.Rows("120:120").Insert() | Out-Null

# Fill in the newly inserted rows.
$ws1.Range("A100").Value = "p99"
$ws1.Range("B100").Value = "p100"
$ws1.Range("C100").Formula = "=A100&B100"

$ws1.Range("A110").Value = "p99.9"
$ws1.Range("B110").Value = "p100"
$ws1.Range("C110").Formula = "=A110&B110"

$ws1.Range("A120").Value = "p99.99"
$ws1.Range("B120").Value = "p100"
$ws1.Range("C120").Formula = "=A120&B120"

# Highlight the four "whole-range to p100" rows in red (100, 110, 120 and the
# pre-existing p99.999-p100 row which is now row 130).
$ws1.Range("A100:C100").Font.Color = 255
$ws1.Range("A110:C110").Font.Color = 255
$ws1.Range("A120:C120").Font.Color = 255
$ws1.Range("A130:C130").Font.Color = 255

$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

$ws1.Range("C1:C130").Select() | Out-Null

# --- Sheet2: insert the matching 3 new columns, mirroring the sheet1 rows ---
$ws2.Columns("CV:CV").Insert() | Out-Null
$ws2.Columns("DF:DF").Insert() | Out-Null
$ws2.Columns("DP:DP").Insert() | Out-Null

$ws2.Range("CV1").Value = "p99p100"
$ws2.Range("DF1").Value = "p99.9p100"
$ws2.Range("DP1").Value = "p99.99p100"

$ws2.Range("CV1").Font.Color = 255
$ws2.Range("DF1").Font.Color = 255
$ws2.Range("DP1").Font.Color = 255
$ws2.Range("DZ1").Font.Color = 255

$ws2.Activate() | Out-Null
$ws2.Range("A1:DZ1").Select() | Out-Null
